$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows imported from the same dataset (IDs 3-5)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Viettel_Digital"
$ws.Range("C4").Value = "192.168.1.12"
$ws.Range("D4").Value = 5000
$ws.Range("E4").Value = $true

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Viettel_CyberSpace"
$ws.Range("C5").Value = "10.10.1.2"
$ws.Range("D5").Value = 3000
$ws.Range("E5").Value = $true

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Viettel_Money"
$ws.Range("C6").Value = "10.10.1.4"
$ws.Range("D6").Value = 4000
$ws.Range("E6").Value = $false

# Widen the data columns to fit the imported content
$ws.Columns.Item(2).ColumnWidth = 24.3
$ws.Columns.Item(3).ColumnWidth = 18.3
$ws.Columns.Item(4).ColumnWidth = 14.85
$ws.Columns.Item(5).ColumnWidth = 17.65

# Leave the selection where the user clicked after importing the rows
$ws.Range("C9").Select() | Out-Null
